# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" stat (previously derived from Strike#). Update the
# computed values for each game row (rows 2-8) to reflect the regenerated
# K values.
$ws.Range("G2").Value = 5
$ws.Range("G3").Value = 7
$ws.Range("G4").Value = 12
$ws.Range("G5").Value = 8
$ws.Range("G6").Value = 7
$ws.Range("G7").Value = 6
$ws.Range("G8").Value = 3
